$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "258.48"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.57%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "26.83"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-2.00%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.666"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "2.40%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05996"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1.82%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.663"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.55%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8582"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.03%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9239"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.67%"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-1.07%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.04805"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "32.79%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07018"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-0.98%"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-5.81%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09131"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001526"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-1.89%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0006056"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.15%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006093"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.09%"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-1.68%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.152"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-1.36%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.168"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-1.52%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3111"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "1.69%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1288"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "0.85%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.131"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "7.41%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04237"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.63%"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.57%"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-5.95%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001200"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.03%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03838"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "0.13%"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "1.22%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.003813"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-38.77%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002419"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "10.03%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01508"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "31.84%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005098"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-6.36%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-0.03%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-30.59%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1303"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-4.17%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.03%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0001998"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.03%"
